# feat: highlight excel update change message
#
# Adds/removes Port-comment (column E) values on the uart / uart_rx / uart_tx
# sheets, and removes the leftover test row (row 16) from the "uart" sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "uart" ---
$wsUart = $wb.Worksheets.Item("uart")

# Drop the stray test row (test_temp / test_port) at the bottom of the table.
$wsUart.Rows.Item(16).Delete()

# Add review comments to a few ports.
$wsUart.Range("E3").Value = "adsfasf"
$wsUart.Range("E5").Value = "adsfa"
$wsUart.Range("E8").Value = "safda"

# --- Sheet "uart_rx" ---
$wsRx = $wb.Worksheets.Item("uart_rx")

$wsRx.Range("E5").Value = "dsaf"
$wsRx.Range("E8").Value = "fgds"
$wsRx.Range("E11").Value = "asdf"

# --- Sheet "uart_tx" ---
$wsTx = $wb.Worksheets.Item("uart_tx")

# These comments were moved elsewhere; clear them out here.
$wsTx.Range("E5").Value = ""
$wsTx.Range("E7").Value = ""
$wsTx.Range("E10").Value = ""
